# Apply the "official names" correction to the District column (G) and
# remove the stray empty Address cells (F) for rows where the address data
# had ended up entirely in the District column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose District (G) column currently reads exactly "Bellary" and
# should become "Ballari (Bellary)"
$bellaryRows = @(3,4,5,7,11,12,13,15,17,19,23,26,27,28,29,30,33,34,38,40)
foreach ($r in $bellaryRows) {
    $ws.Cells.Item($r, 7).Value = "Ballari (Bellary)"
}

# Rows whose District (G) column currently reads exactly "Madhugiri" or
# "Tumkur" and should become "Tumakuru (Tumkur)"
$tumakuruRows = @(8,16,18,20,22,24,25,32,39)
foreach ($r in $tumakuruRows) {
    $ws.Cells.Item($r, 7).Value = "Tumakuru (Tumkur)"
}

# Remove the stray empty inline-string cells in column F (Address) for the
# rows where the address text had been mistakenly placed in G instead.
$emptyFRows = @(9,10,14,31,36)
foreach ($r in $emptyFRows) {
    $ws.Cells.Item($r, 6).ClearContents()
}
